# Applies the edits described in the commit/diff:
#  - Productdata sheet: reduce StartingInventories (col C) and BackorderCosts
#    (col E) values for rows 2-9 (average slowmoving related tuning), and
#    tweak G7:G9.
#  - Capacity sheet: rescale column B (capacity) values.
#  - ProcessingTime sheet: adjust a handful of single cells.

$wb = $excel.ActiveWorkbook

# --- Productdata sheet ---
$ws = $wb.Worksheets.Item("Productdata")

$ws.Cells.Item(2, 3).Value = 2        # C2: 11 -> 2
$ws.Cells.Item(2, 5).Value = 0.41     # E2: 2.36 -> 0.41

$ws.Cells.Item(3, 3).Value = 3        # C3: 17 -> 3
$ws.Cells.Item(3, 5).Value = 0.47     # E3: 2.1 -> 0.47

$ws.Cells.Item(4, 3).Value = 4        # C4: 21 -> 4
$ws.Cells.Item(4, 5).Value = 0.86     # E4: 3.45 -> 0.86

$ws.Cells.Item(5, 3).Value = 6        # C5: 35 -> 6
$ws.Cells.Item(5, 5).Value = 0.72     # E5: 3.37 -> 0.72

$ws.Cells.Item(6, 3).Value = 4        # C6: 37 -> 4
$ws.Cells.Item(6, 5).Value = 0.91     # E6: 3.81 -> 0.91

$ws.Cells.Item(7, 3).Value = 1        # C7: 5 -> 1
$ws.Cells.Item(7, 5).Value = 0.21     # E7: 1.07 -> 0.21
$ws.Cells.Item(7, 7).Value = 1        # G7: 5 -> 1

$ws.Cells.Item(8, 3).Value = 1        # C8: 5 -> 1
$ws.Cells.Item(8, 5).Value = 0.22     # E8: 0.82 -> 0.22
$ws.Cells.Item(8, 7).Value = 1        # G8: 5 -> 1

$ws.Cells.Item(9, 3).Value = 1        # C9: 5 -> 1
$ws.Cells.Item(9, 5).Value = 0.2      # E9: 1.18 -> 0.2
$ws.Cells.Item(9, 7).Value = 1        # G9: 5 -> 1

# --- Capacity sheet ---
$ws = $wb.Worksheets.Item("Capacity")

$ws.Cells.Item(2, 2).Value = 14.4     # B2: 36 -> 14.4
$ws.Cells.Item(3, 2).Value = 3.6      # B3: 18 -> 3.6
$ws.Cells.Item(4, 2).Value = 21.6     # B4: 180 -> 21.6
$ws.Cells.Item(5, 2).Value = 28.8     # B5: 36 -> 28.8
$ws.Cells.Item(6, 2).Value = 21.6     # B6: 144 -> 21.6
$ws.Cells.Item(8, 2).Value = 9        # B8: 18 -> 9
$ws.Cells.Item(9, 2).Value = 7.2      # B9: 45 -> 7.2

# --- ProcessingTime sheet ---
$ws = $wb.Worksheets.Item("ProcessingTime")

$ws.Cells.Item(2, 2).Value = 4        # B2: 2 -> 4
$ws.Cells.Item(4, 4).Value = 3        # D4: 5 -> 3
$ws.Cells.Item(5, 5).Value = 4        # E5: 1 -> 4
$ws.Cells.Item(6, 6).Value = 3        # F6: 4 -> 3
$ws.Cells.Item(7, 7).Value = 5        # G7: 1 -> 5
$ws.Cells.Item(8, 8).Value = 5        # H8: 2 -> 5
$ws.Cells.Item(9, 9).Value = 4        # I9: 5 -> 4
